$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.989.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.91%  "
$ws.Range("D3").Value = "'1.600.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.61%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'211.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.67%  "
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("E7").Value = "  +1.55%  "
$ws.Range("D8").Value = "'0.247"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.93%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").Value = "'18.10"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.00%  "
$ws.Range("D11").Value = "'0.0811"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.62%  "
$ws.Range("D12").Value = "'1.824.28"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.72%  "
$ws.Range("D13").Value = "'1.602.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.47%  "
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("D16").Value = "'26.002.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.91%  "
$ws.Range("D17").Value = "'60.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.61%  "
$ws.Range("D18").Value = "'0.0₃0722"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.36%  "
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("D20").Value = "'201.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.16%  "
$ws.Range("D21").Value = "'4.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.79%  "
$ws.Range("D22").Value = "'9.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("E23").Value = "  +2.48%  "
$ws.Range("D24").Value = "'1.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.73%  "
$ws.Range("D25").Value = "'141.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("E27").Value = "  -7.04%  "
$ws.Range("D28").Value = "'15.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.34%  "
$ws.Range("E29").Value = "  +0.73%  "
$ws.Range("E30").Value = "  +2.14%  "
$ws.Range("E31").Value = "  +2.23%  "
$ws.Range("D32").Value = "'3.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.57%  "
$ws.Range("D33").Value = "'2.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.06%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  +1.99%  "
$ws.Range("D36").Value = "'1.125.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.86%  "
$ws.Range("D37").Value = "'0.0163"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.84%  "
$ws.Range("E38").Value = "  -0.25%  "
$ws.Range("D39").Value = "'0.790"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.56%  "
$ws.Range("E40").Value = "  -0.47%  "
$ws.Range("E41").Value = "  -0.81%  "
$ws.Range("E42").Value = "  -1.36%  "
$ws.Range("B43").Value = "RocketPoolETH"
$ws.Range("C43").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D43").Value = "'1.735.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.61%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'5.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.20%  "
$ws.Range("D45").Value = "'92.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("D46").Value = "'1.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.92%  "
$ws.Range("D47").Value = "'53.57"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.08%  "
$ws.Range("D48").Value = "'0.0503"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.20%  "
$ws.Range("D49").Value = "'0.408"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.06%  "
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.58%  "

Write-Output "Applied 86 cell updates"
